# Updated cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'62.142.74"
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = "'2.437.95"
$ws.Range('E3').Value = '  +4.97%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'556.55"
$ws.Range('E5').Value = '  +2.12%  '
$ws.Range('D6').Value = "'138.78"
$ws.Range('E6').Value = '  +5.90%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'0.586"
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('D9').Value = "'2.435.64"
$ws.Range('E9').Value = '  +4.98%  '
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('D11').Value = "'5.75"
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = "'0.348"
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').Value = "'26.02"
$ws.Range('E14').Value = '  +10.48%  '
$ws.Range('D15').Value = "'2.871.65"
$ws.Range('E15').Value = '  +4.97%  '
$ws.Range('D16').Value = "'62.045.66"
$ws.Range('E16').Value = '  +3.11%  '
$ws.Range('E17').Value = '  +5.92%  '
$ws.Range('D18').Value = "'2.437.68"
$ws.Range('E18').Value = '  +5.15%  '
$ws.Range('D19').Value = "'11.17"
$ws.Range('E19').Value = '  +5.69%  '
$ws.Range('D20').Value = "'346.05"
$ws.Range('E20').Value = '  +10.12%  '
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('D22').Value = "'6.83"
$ws.Range('E22').Value = '  +3.19%  '
$ws.Range('D24').Value = "'65.29"
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E27').Value = '  +11.27%  '
$ws.Range('D28').Value = "'8.27"
$ws.Range('E28').Value = '  +6.00%  '
$ws.Range('D29').Value = "'1.35"
$ws.Range('E29').Value = '  +12.53%  '
$ws.Range('D30').Value = "'0.0₃0787"
$ws.Range('E30').Value = '  +7.67%  '
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').Value = "'6.35"
$ws.Range('E32').Value = '  +7.03%  '
$ws.Range('D33').Value = "'171.37"
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').Value = "'1.47"
$ws.Range('E34').Value = '  +6.80%  '
$ws.Range('D35').Value = "'0.398"
$ws.Range('E35').Value = '  +4.94%  '
$ws.Range('D36').Value = "'18.62"
$ws.Range('E36').Value = '  +4.69%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'4.51"
$ws.Range('E37').Value = '  +11.60%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = "'368.81"
$ws.Range('E38').Value = '  +15.72%  '
$ws.Range('D39').Value = "'0.998"
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +11.00%  '
$ws.Range('D42').Value = "'39.33"
$ws.Range('E42').Value = '  +3.86%  '
$ws.Range('D43').Value = "'146.91"
$ws.Range('E43').Value = '  +6.88%  '
$ws.Range('D44').Value = "'3.69"
$ws.Range('E44').Value = '  +6.24%  '
$ws.Range('D45').Value = "'20.70"
$ws.Range('E45').Value = '  +8.94%  '
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('D47').Value = "'0.589"
$ws.Range('E47').Value = '  +4.77%  '
$ws.Range('E48').Value = '  +5.25%  '
$ws.Range('D49').Value = "'18.05"
$ws.Range('E49').Value = '  +7.07%  '
$ws.Range('E50').Value = '  +4.28%  '
$ws.Range('D51').Value = "'1.72"
$ws.Range('E51').Value = '  +12.49%  '
